# Workbook: MassWateR example WQX template (sheet "Meta")
# Change: the default "Collection/Monitoring Organization" value ("MassWateR")
# used for TP, TDP and E.coli rows is renamed to "MassBays".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meta")

$ws.Range("B4").Value = "MassBays"
$ws.Range("B5").Value = "MassBays"
$ws.Range("B6").Value = "MassBays"

# Update the active cell/selection to reflect where the edit left off.
$ws.Range("F9").Select() | Out-Null
